$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing last row (178) down to make room for two new rows.
# First insert a row at 178 (old row 178 -> 179), then insert a row at 180
# (after the shifted row) for the new trailing row.
$ws.Rows("178:178").Insert()

# New row 178 (inserted before the former last row)
$ws.Range("A178").Value = 45454.2916666667
$ws.Range("B178").Value = 78500
$ws.Range("C178").Value = 3.40000009536743
$ws.Range("D178").Value = 3.15000009536743
$ws.Range("E178").Value = 3.23000001907349
$ws.Range("F178").Value = 3.15000009536743
$ws.Range("G178").NumberFormat = "@"
$ws.Range("G178").Value = "3.15000009536743"
$ws.Range("G178").ClearFormats()
$ws.Range("H178").Value = "EAV.MI"

# Row 179 is the former row 178, shifted down; only the date changed.
$ws.Range("A179").Value = 45455.2916666667

# New row 180, appended after the shifted row.
$ws.Range("A179").Copy()
$ws.Range("A180").PasteSpecial(-4122)
$ws.Range("A180").Value = 45456.5604282407
$ws.Range("B180").Value = 14000
$ws.Range("C180").Value = 3.22000002861023
$ws.Range("D180").Value = 3.15000009536743
$ws.Range("E180").Value = 3.20000004768372
$ws.Range("F180").Value = 3.22000002861023
$ws.Range("G180").NumberFormat = "@"
$ws.Range("G180").Value = "3.22000002861023"
$ws.Range("G180").ClearFormats()
$ws.Range("H180").Value = "EAV.MI"
